$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.519.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.81%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.628.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.05%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'202.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +7.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'580.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.80%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.621.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.98%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.624"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +3.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.19%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.687"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.38%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'61.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +17.04%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.151"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.78%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +13.67%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'10.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +4.96%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.202.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.96%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.635.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.07%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  +1.09%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'19.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +6.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +4.98%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'68.414.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.01%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +3.42%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'406.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.87%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +18.62%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'4.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.81%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'86.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.06%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "'Toncoin"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'4.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +15.61%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "'ImmutableX"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.10%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'12.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.81%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +2.32%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'9.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +8.47%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +11.12%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'31.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.47%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'677.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +9.09%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'12.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.00%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +3.46%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'63.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.36%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'42.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.87%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.421"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +8.19%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.16%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.0₃0779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +5.51%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'3.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +17.38%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.210.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +8.77%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.135"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.96%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.12%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Fetch.AI"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'2.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +11.75%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +26.95%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +16.09%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +5.40%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Stellar"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.132"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +2.65%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'THORChain"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'8.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +6.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'3.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.16%  "
$ws.Range("E51").Style = "Normal"
